$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new row at 131, copying formatting from the row above (row 130)
# so fills/number formats/wrap text carry over exactly like a manual
# "Insert Copied Cells" in Excel.
$ws.Rows(130).Copy()
$ws.Rows(131).Insert()

# Populate the new row's contents (acronym entry for "indst").
$ws.Cells.Item(131, 1).Value = "indst"
$ws.Cells.Item(131, 2).Value = "ItICM"
$ws.Cells.Item(131, 3).Value = "Industry to ISIC Code Map"
$ws.Cells.Item(131, 4).ClearContents()
$ws.Cells.Item(131, 5).ClearContents()
$ws.Cells.Item(131, 6).Value = "medium"
